$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B48 currently holds the text "3" stored as an inline string; convert it to a
# true numeric value (3) while leaving everything else in the row untouched.
$ws.Range("B48").Value = 3

# Append new row 49 with the additional annotation record.
$ws.Range("A49").Value = "Ruilin"

# B49 must stay a text value "4" (not numeric), so force text formatting
# before assigning, then restore the default "Normal" style so no stray
# number-format style gets attached to the cell.
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "4"
$ws.Range("B49").Style = "Normal"

$ws.Range("C49").Value = "thank"
$ws.Range("D49").Value = "ACK"
$ws.Range("E49").Value = "OTH"
$ws.Range("F49").Value = "658343d9-2c6f-4c77-9518-16756d4b8755"
$ws.Range("G49").Value = "SylJ1D1C-_annotated.xlsx"
$ws.Range("H49").Value = "First, we would like to thank the reviewer for carefully evaluating our paper."
